$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3987
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = 229
$ws.Cells.Item(2, 5).Value = 535
$ws.Cells.Item(2, 6).Value = 15
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 891
$ws.Cells.Item(2, 9).Value = 31702
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(2, 11).Value = 283
$ws.Cells.Item(2, 12).Value = 15833
$ws.Cells.Item(3, 2).Value = 821
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 821
$ws.Cells.Item(3, 5).Value = 822
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 9).Value = 8786
$ws.Cells.Item(3, 10).Value = 17
$ws.Cells.Item(4, 2).Value = 258
$ws.Cells.Item(5, 2).Value = 803
$ws.Cells.Item(5, 4).Value = 336
$ws.Cells.Item(5, 5).Value = 336
$ws.Cells.Item(5, 9).Value = 9014
$ws.Cells.Item(5, 10).Value = 9
$ws.Cells.Item(6, 2).Value = 259
$ws.Cells.Item(7, 2).Value = 255
$ws.Cells.Item(7, 4).Value = 13
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 9).Value = 6154
$ws.Cells.Item(8, 2).Value = 772
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 302
$ws.Cells.Item(8, 5).Value = 329
$ws.Cells.Item(8, 6).Value = 20
$ws.Cells.Item(8, 8).Value = 621
$ws.Cells.Item(8, 9).Value = 3775
$ws.Cells.Item(9, 2).Value = 774
$ws.Cells.Item(9, 4).Value = 73
$ws.Cells.Item(9, 5).Value = 73
$ws.Cells.Item(9, 9).Value = 9014
$ws.Cells.Item(10, 2).Value = 576
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 63
$ws.Cells.Item(10, 5).Value = 72
$ws.Cells.Item(10, 6).Value = 7
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 1429
$ws.Cells.Item(10, 9).Value = 18772
$ws.Cells.Item(11, 2).Value = 754
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 300
$ws.Cells.Item(11, 5).Value = 322
$ws.Cells.Item(11, 6).Value = 19
$ws.Cells.Item(11, 8).Value = 2967
$ws.Cells.Item(11, 9).Value = 32145
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(12, 2).Value = 427
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 27
$ws.Cells.Item(12, 5).Value = 29
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 8).Value = 556
$ws.Cells.Item(12, 9).Value = 15882
$ws.Cells.Item(13, 2).Value = 461
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = 325
$ws.Cells.Item(13, 5).Value = 378
$ws.Cells.Item(13, 6).Value = 47
$ws.Cells.Item(13, 7).Value = 3
$ws.Cells.Item(13, 8).Value = 3963
$ws.Cells.Item(13, 9).Value = 6822
$ws.Cells.Item(13, 10).Value = 9
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 33
$ws.Cells.Item(14, 2).Value = 804
$ws.Cells.Item(14, 4).Value = 273
$ws.Cells.Item(14, 5).Value = 275
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 8).Value = 36
$ws.Cells.Item(14, 9).Value = 9118
$ws.Cells.Item(14, 10).Value = 8
$ws.Cells.Item(15, 2).Value = 257
$ws.Cells.Item(16, 2).Value = 67
$ws.Cells.Item(16, 3).Value = 13
$ws.Cells.Item(16, 4).Value = 518
$ws.Cells.Item(16, 5).Value = 617
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 8).Value = 776
$ws.Cells.Item(16, 9).Value = 15334
$ws.Cells.Item(16, 10).Value = 6
$ws.Cells.Item(16, 11).Value = 70
$ws.Cells.Item(16, 12).Value = 1528
$ws.Cells.Item(17, 2).Value = 57
$ws.Cells.Item(17, 4).Value = 28
$ws.Cells.Item(17, 5).Value = 28
$ws.Cells.Item(17, 9).Value = 7143
$ws.Cells.Item(18, 2).Value = 5
$ws.Cells.Item(18, 4).Value = 516
$ws.Cells.Item(18, 5).Value = 554
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(18, 8).Value = 72
$ws.Cells.Item(18, 9).Value = 7665
$ws.Cells.Item(18, 10).Value = 10
$ws.Cells.Item(18, 11).Value = 33
$ws.Cells.Item(18, 12).Value = 596
$ws.Cells.Item(19, 2).Value = 771
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 173
$ws.Cells.Item(19, 5).Value = 177
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(19, 8).Value = 197
$ws.Cells.Item(19, 9).Value = 9338
$ws.Cells.Item(20, 2).Value = 1383
$ws.Cells.Item(21, 2).Value = 364
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 4).Value = 78
$ws.Cells.Item(21, 5).Value = 80
$ws.Cells.Item(21, 9).Value = 17700
$ws.Cells.Item(22, 2).Value = 181
$ws.Cells.Item(22, 4).Value = 3
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 9).Value = 3333
$ws.Cells.Item(23, 2).Value = 813
$ws.Cells.Item(23, 4).Value = 21
$ws.Cells.Item(23, 5).Value = 21
$ws.Cells.Item(23, 9).Value = 6667
$ws.Cells.Item(23, 10).Value = 4
$ws.Cells.Item(24, 2).Value = 897
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 26
$ws.Cells.Item(24, 5).Value = 27
$ws.Cells.Item(24, 9).Value = 9878
$ws.Cells.Item(25, 2).Value = 69
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 61
$ws.Cells.Item(25, 5).Value = 65
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 156
$ws.Cells.Item(25, 9).Value = 5082
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 12).Value = 308
$ws.Cells.Item(26, 2).Value = 124
$ws.Cells.Item(26, 4).Value = 44
$ws.Cells.Item(26, 5).Value = 46
$ws.Cells.Item(26, 8).Value = 435
$ws.Cells.Item(26, 9).Value = 7442
$ws.Cells.Item(26, 10).Value = 1
$ws.Cells.Item(27, 2).Value = 4
$ws.Cells.Item(27, 3).Value = 4
$ws.Cells.Item(27, 4).Value = 99
$ws.Cells.Item(27, 5).Value = 109
$ws.Cells.Item(27, 6).Value = 6
$ws.Cells.Item(27, 8).Value = 550
$ws.Cells.Item(27, 9).Value = 4242
